$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before the current row 6 ("ready"/"go" get pushed down
# to rows 13-14), making room for the new HUD / options strings.
$ws.Rows("6:12").Insert()

# Populate the new rows. The shared-strings table records the order the
# values were first entered, so we write them in that original authoring
# order (options, music, sound, on, off, speech, close) even though the
# rows end up in a different order on the sheet.
$ws.Cells.Item(6, 1).Value2 = "options"
$ws.Cells.Item(6, 2).Value2 = "OPTIONS"

$ws.Cells.Item(7, 1).Value2 = "music"
$ws.Cells.Item(7, 2).Value2 = "MUSIC"

$ws.Cells.Item(8, 1).Value2 = "sound"
$ws.Cells.Item(8, 2).Value2 = "SOUND"

$ws.Cells.Item(11, 1).Value2 = "on"
$ws.Cells.Item(11, 2).Value2 = "ON"

$ws.Cells.Item(12, 1).Value2 = "off"
$ws.Cells.Item(12, 2).Value2 = "OFF"

$ws.Cells.Item(9, 1).Value2 = "speech"
$ws.Cells.Item(9, 2).Value2 = "SPEECH"

$ws.Cells.Item(10, 1).Value2 = "close"
$ws.Cells.Item(10, 2).Value2 = "CLOSE"

# All of the new Value cells (B6:B12) wrap text, matching the style used
# by the other descriptive rows on the sheet.
$ws.Range("B6:B12").WrapText = $true

# Update the active selection to match the edited document state.
$ws.Range("B11").Select() | Out-Null
